# Auto-generated Excel COM-interop script to apply the Excalibur_Profits diff
# Updates computed market-price / profit columns (H-N) across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 10).Value = 363.66666  # J19
$ws.Cells.Item(19, 11).Value = 894.875  # K19
$ws.Cells.Item(19, 13).Value = -719.875  # M19
$ws.Cells.Item(19, 14).Value = -713.66666  # N19
$ws.Cells.Item(19, 12).Value = 363.66666  # L19
$ws.Cells.Item(19, 8).Value = 667.2143  # H19
$ws.Cells.Item(19, 9).Value = 894.875  # I19
$ws.Cells.Item(25, 8).Value = 25000  # H25
$ws.Cells.Item(25, 14).Value = -75348  # N25
$ws.Cells.Item(25, 12).Value = 75000  # L25
$ws.Cells.Item(25, 10).Value = 25000  # J25
$ws.Cells.Item(33, 13).Value = -63  # M33
$ws.Cells.Item(33, 8).Value = 536.3333  # H33
$ws.Cells.Item(33, 9).Value = 292  # I33
$ws.Cells.Item(33, 11).Value = 292  # K33
$ws.Cells.Item(98, 14).Value = -3986.3333  # N98
$ws.Cells.Item(98, 8).Value = 1705.8077  # H98
$ws.Cells.Item(98, 12).Value = 990.3333  # L98
$ws.Cells.Item(98, 9).Value = 1799.1305  # I98
$ws.Cells.Item(98, 13).Value = -301.1305  # M98
$ws.Cells.Item(98, 10).Value = 990.3333  # J98
$ws.Cells.Item(98, 11).Value = 1799.1305  # K98
$ws.Cells.Item(116, 13).Value = -7385.031000000001  # M116
$ws.Cells.Item(116, 9).Value = 10827.031  # I116
$ws.Cells.Item(116, 11).Value = 10827.031  # K116
$ws.Cells.Item(116, 8).Value = 11995.565  # H116
$ws.Cells.Item(122, 14).Value = -7870.9999  # N122
$ws.Cells.Item(122, 12).Value = 2970.9999  # L122
$ws.Cells.Item(122, 10).Value = 990.3333  # J122
$ws.Cells.Item(122, 13).Value = -2947.3915  # M122
$ws.Cells.Item(122, 8).Value = 1705.8077  # H122
$ws.Cells.Item(122, 11).Value = 5397.3915  # K122
$ws.Cells.Item(122, 9).Value = 1799.1305  # I122
$ws.Cells.Item(129, 8).Value = 1894.9231  # H129
$ws.Cells.Item(129, 13).Value = -230.7999999999993  # M129
$ws.Cells.Item(129, 11).Value = 5230.799999999999  # K129
$ws.Cells.Item(129, 9).Value = 1743.6  # I129
$ws.Cells.Item(129, 14).Value = -17197.9999  # N129
$ws.Cells.Item(129, 10).Value = 2399.3333  # J129
$ws.Cells.Item(129, 12).Value = 7197.999899999999  # L129
$ws.Cells.Item(132, 11).Value = 117779.085  # K132
$ws.Cells.Item(132, 9).Value = 39259.695  # I132
$ws.Cells.Item(132, 8).Value = 35978.363  # H132
$ws.Cells.Item(132, 13).Value = -115249.085  # M132
$ws.Cells.Item(138, 10).Value = 3648.2979  # J138
$ws.Cells.Item(138, 8).Value = 4060.1228  # H138
$ws.Cells.Item(138, 12).Value = 10944.8937  # L138
$ws.Cells.Item(138, 14).Value = -21224.8937  # N138
$ws.Cells.Item(141, 8).Value = 1800  # H141
$ws.Cells.Item(141, 13).Value = -220  # M141
$ws.Cells.Item(141, 9).Value = 1800  # I141
$ws.Cells.Item(141, 11).Value = 5400  # K141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 11).Value = 3066.2126  # K32
$ws.Cells.Item(32, 13).Value = -2779.2126  # M32
$ws.Cells.Item(32, 14).Value = -14573.667  # N32
$ws.Cells.Item(32, 8).Value = 3722.22  # H32
$ws.Cells.Item(32, 9).Value = 3066.2126  # I32
$ws.Cells.Item(32, 12).Value = 13999.667  # L32
$ws.Cells.Item(32, 10).Value = 13999.667  # J32
$ws.Cells.Item(74, 11).Value = 3679489.8  # K74
$ws.Cells.Item(74, 13).Value = -3678615.8  # M74
$ws.Cells.Item(74, 9).Value = 3679489.8  # I74
$ws.Cells.Item(74, 8).Value = 2979637.2  # H74
$ws.Cells.Item(77, 8).Value = 2979637.2  # H77
$ws.Cells.Item(77, 13).Value = -18393081  # M77
$ws.Cells.Item(77, 9).Value = 3679489.8  # I77
$ws.Cells.Item(77, 11).Value = 18397449  # K77
$ws.Cells.Item(110, 13).Value = 484.9000000000001  # M110
$ws.Cells.Item(110, 9).Value = 1560.1  # I110
$ws.Cells.Item(110, 11).Value = 1560.1  # K110
$ws.Cells.Item(110, 8).Value = 1718.9166  # H110
$ws.Cells.Item(132, 10).Value = 6088.647  # J132
$ws.Cells.Item(132, 11).Value = 1523502.48  # K132
$ws.Cells.Item(132, 9).Value = 507834.16  # I132
$ws.Cells.Item(132, 8).Value = 372442.5  # H132
$ws.Cells.Item(132, 14).Value = -23325.941  # N132
$ws.Cells.Item(132, 13).Value = -1520972.48  # M132
$ws.Cells.Item(132, 12).Value = 18265.941  # L132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2284.1072  # H20
$ws.Cells.Item(20, 13).Value = -1712.9048  # M20
$ws.Cells.Item(20, 12).Value = 3256.7144  # L20
$ws.Cells.Item(20, 9).Value = 1959.9048  # I20
$ws.Cells.Item(20, 10).Value = 3256.7144  # J20
$ws.Cells.Item(20, 14).Value = -3750.7144  # N20
$ws.Cells.Item(20, 11).Value = 1959.9048  # K20
$ws.Cells.Item(99, 9).Value = 1014.75  # I99
$ws.Cells.Item(99, 11).Value = 1014.75  # K99
$ws.Cells.Item(99, 13).Value = 483.25  # M99
$ws.Cells.Item(99, 8).Value = 2312.7  # H99
$ws.Cells.Item(105, 11).Value = 2589.4  # K105
$ws.Cells.Item(105, 9).Value = 2589.4  # I105
$ws.Cells.Item(105, 13).Value = -842.4000000000001  # M105
$ws.Cells.Item(105, 8).Value = 3310.3809  # H105
$ws.Cells.Item(134, 9).Value = 683409.0600000001  # I134
$ws.Cells.Item(134, 11).Value = 2050227.18  # K134
$ws.Cells.Item(134, 8).Value = 499687.75  # H134
$ws.Cells.Item(134, 13).Value = -2047692.18  # M134
$ws.Cells.Item(139, 12).Value = 55000  # L139
$ws.Cells.Item(139, 8).Value = 55000  # H139
$ws.Cells.Item(139, 10).Value = 55000  # J139
$ws.Cells.Item(139, 14).Value = -65280  # N139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 9).Value = 3137.0908  # I31
$ws.Cells.Item(31, 8).Value = 9037.200000000001  # H31
$ws.Cells.Item(31, 11).Value = 3137.0908  # K31
$ws.Cells.Item(31, 13).Value = -2842.0908  # M31
$ws.Cells.Item(34, 8).Value = 9037.200000000001  # H34
$ws.Cells.Item(34, 13).Value = -2935.0908  # M34
$ws.Cells.Item(34, 11).Value = 3137.0908  # K34
$ws.Cells.Item(34, 9).Value = 3137.0908  # I34
$ws.Cells.Item(107, 8).Value = 1115.1666  # H107
$ws.Cells.Item(107, 9).Value = 969.0833  # I107
$ws.Cells.Item(107, 13).Value = 950.9167  # M107
$ws.Cells.Item(107, 11).Value = 969.0833  # K107
$ws.Cells.Item(122, 13).Value = -160.6666  # M122
$ws.Cells.Item(122, 8).Value = 893.8182  # H122
$ws.Cells.Item(122, 11).Value = 2610.6666  # K122
$ws.Cells.Item(122, 9).Value = 870.2222  # I122
$ws.Cells.Item(134, 9).Value = 1856.6389  # I134
$ws.Cells.Item(134, 11).Value = 5569.9167  # K134
$ws.Cells.Item(134, 8).Value = 2695.825  # H134
$ws.Cells.Item(134, 13).Value = -3034.9167  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 5347.4287  # H3
$ws.Cells.Item(3, 9).Value = 1266.909  # I3
$ws.Cells.Item(3, 11).Value = 3800.727  # K3
$ws.Cells.Item(3, 13).Value = -3688.727  # M3
$ws.Cells.Item(96, 8).Value = 6439.3335  # H96
$ws.Cells.Item(96, 10).Value = 7322.2  # J96
$ws.Cells.Item(96, 12).Value = 21966.6  # L96
$ws.Cells.Item(96, 14).Value = -26084.6  # N96
$ws.Cells.Item(114, 9).Value = 578.8333  # I114
$ws.Cells.Item(114, 10).Value = 2623.6875  # J114
$ws.Cells.Item(114, 14).Value = -14379.0625  # N114
$ws.Cells.Item(114, 12).Value = 7871.0625  # L114
$ws.Cells.Item(114, 8).Value = 1747.3214  # H114
$ws.Cells.Item(114, 13).Value = 1517.5001  # M114
$ws.Cells.Item(114, 11).Value = 1736.4999  # K114
$ws.Cells.Item(131, 12).Value = 53210.571  # L131
$ws.Cells.Item(131, 14).Value = -63290.571  # N131
$ws.Cells.Item(131, 10).Value = 17736.857  # J131
$ws.Cells.Item(131, 8).Value = 16832.25  # H131
$ws.Cells.Item(140, 9).Value = 1117.5555  # I140
$ws.Cells.Item(140, 11).Value = 3352.6665  # K140
$ws.Cells.Item(140, 8).Value = 1187.3243  # H140
$ws.Cells.Item(140, 13).Value = 1827.3335  # M140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 9).Value = 5200000  # I14
$ws.Cells.Item(14, 11).Value = 5200000  # K14
$ws.Cells.Item(14, 13).Value = -5199832  # M14
$ws.Cells.Item(14, 8).Value = 14466667  # H14
$ws.Cells.Item(15, 10).Value = 40750  # J15
$ws.Cells.Item(15, 11).Value = 4200  # K15
$ws.Cells.Item(15, 8).Value = 33440  # H15
$ws.Cells.Item(15, 13).Value = -3912  # M15
$ws.Cells.Item(15, 14).Value = -41326  # N15
$ws.Cells.Item(15, 12).Value = 40750  # L15
$ws.Cells.Item(15, 9).Value = 4200  # I15
$ws.Cells.Item(81, 9).Value = 4200  # I81
$ws.Cells.Item(81, 8).Value = 33440  # H81
$ws.Cells.Item(81, 12).Value = 40750  # L81
$ws.Cells.Item(81, 10).Value = 40750  # J81
$ws.Cells.Item(81, 14).Value = -42746  # N81
$ws.Cells.Item(81, 11).Value = 4200  # K81
$ws.Cells.Item(81, 13).Value = -3202  # M81
$ws.Cells.Item(84, 9).Value = 4200  # I84
$ws.Cells.Item(84, 10).Value = 40750  # J84
$ws.Cells.Item(84, 13).Value = -7608  # M84
$ws.Cells.Item(84, 8).Value = 33440  # H84
$ws.Cells.Item(84, 14).Value = -132234  # N84
$ws.Cells.Item(84, 12).Value = 122250  # L84
$ws.Cells.Item(84, 11).Value = 12600  # K84
$ws.Cells.Item(97, 14).Value = -3107  # N97
$ws.Cells.Item(97, 8).Value = 1091.6666  # H97
$ws.Cells.Item(97, 12).Value = 2115  # L97
$ws.Cells.Item(97, 10).Value = 2115  # J97
$ws.Cells.Item(113, 9).Value = 1572.5834  # I113
$ws.Cells.Item(113, 8).Value = 2868.55  # H113
$ws.Cells.Item(113, 13).Value = 597.4166  # M113
$ws.Cells.Item(113, 11).Value = 1572.5834  # K113
$ws.Cells.Item(122, 14).Value = -50893.999  # N122
$ws.Cells.Item(122, 12).Value = 45993.999  # L122
$ws.Cells.Item(122, 10).Value = 15331.333  # J122
$ws.Cells.Item(122, 13).Value = -2099.6764  # M122
$ws.Cells.Item(122, 8).Value = 2636.6758  # H122
$ws.Cells.Item(122, 11).Value = 4549.6764  # K122
$ws.Cells.Item(122, 9).Value = 1516.5588  # I122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 12).Value = 2266  # L93
$ws.Cells.Item(93, 10).Value = 2266  # J93
$ws.Cells.Item(93, 14).Value = -4762  # N93
$ws.Cells.Item(93, 8).Value = 2219.6  # H93
$ws.Cells.Item(100, 11).Value = 2150.4285  # K100
$ws.Cells.Item(100, 9).Value = 2150.4285  # I100
$ws.Cells.Item(100, 8).Value = 8059.5  # H100
$ws.Cells.Item(100, 13).Value = -1609.4285  # M100
$ws.Cells.Item(132, 10).Value = 5879.9  # J132
$ws.Cells.Item(132, 11).Value = 3152500.5  # K132
$ws.Cells.Item(132, 9).Value = 1050833.5  # I132
$ws.Cells.Item(132, 8).Value = 807821  # H132
$ws.Cells.Item(132, 14).Value = -22699.7  # N132
$ws.Cells.Item(132, 13).Value = -3149970.5  # M132
$ws.Cells.Item(132, 12).Value = 17639.7  # L132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 10).Value = 0  # J19
$ws.Cells.Item(19, 14).ClearContents()  # N19
$ws.Cells.Item(19, 12).Value = 0  # L19
$ws.Cells.Item(19, 8).Value = 0  # H19
$ws.Cells.Item(105, 12).Value = 70615  # L105
$ws.Cells.Item(105, 10).Value = 70615  # J105
$ws.Cells.Item(105, 8).Value = 70615  # H105
$ws.Cells.Item(105, 14).Value = -77603  # N105
$ws.Cells.Item(122, 13).Value = -2274.090999999999  # M122
$ws.Cells.Item(122, 8).Value = 2074.125  # H122
$ws.Cells.Item(122, 11).Value = 4724.090999999999  # K122
$ws.Cells.Item(122, 9).Value = 1574.697  # I122
$ws.Cells.Item(123, 10).Value = 47000  # J123
$ws.Cells.Item(123, 8).Value = 47000  # H123
$ws.Cells.Item(123, 12).Value = 47000  # L123
$ws.Cells.Item(123, 14).Value = -56800  # N123
$ws.Cells.Item(132, 11).Value = 16324773  # K132
$ws.Cells.Item(132, 9).Value = 5441591  # I132
$ws.Cells.Item(132, 8).Value = 3533526.5  # H132
$ws.Cells.Item(132, 13).Value = -16322243  # M132
